$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = "62.922.46"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "3.063.31"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "537.04"
$ws.Range("E5").Value = "  -3.19%  "

Set-TextValue $ws.Range("D6") "133.30"
$ws.Range("E6").Value = "  -2.76%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.058.52"
$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("E10").Value = "  -3.39%  "

Set-TextValue $ws.Range("D11") "6.07"
$ws.Range("E11").Value = "  -8.09%  "

Set-TextValue $ws.Range("D12") "0.452"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("E13").Value = "  +2.86%  "

Set-TextValue $ws.Range("D14") "34.22"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").Value = "3.558.25"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "62.888.20"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").Value = "3.063.29"
$ws.Range("E18").Value = "  -0.90%  "

Set-TextValue $ws.Range("D19") "6.63"
$ws.Range("E19").Value = "  -0.41%  "

Set-TextValue $ws.Range("D20") "482.31"
$ws.Range("E20").Value = "  -3.55%  "

Set-TextValue $ws.Range("D21") "13.32"
$ws.Range("E21").Value = "  -1.47%  "

Set-TextValue $ws.Range("D22") "0.695"
$ws.Range("E22").Value = "  -1.58%  "

Set-TextValue $ws.Range("D23") "7.10"
$ws.Range("E23").Value = "  -2.33%  "

Set-TextValue $ws.Range("D24") "79.11"
$ws.Range("E24").Value = "  +1.42%  "

Set-TextValue $ws.Range("D25") "12.09"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("E26").Value = "  +0.01%  "

Set-TextValue $ws.Range("D27") "2.70"
$ws.Range("E27").Value = "  -2.23%  "

Set-TextValue $ws.Range("D28") "8.09"
$ws.Range("E28").Value = "  -0.89%  "

Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.06%  "

Set-TextValue $ws.Range("D30") "25.95"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("E31").Value = "  -7.27%  "

$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("E33").Value = "  -6.26%  "

Set-TextValue $ws.Range("D34") "57.10"
$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D35") "5.33"
$ws.Range("E35").Value = "  +3.62%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D36") "6.01"
$ws.Range("E36").Value = "  +2.32%  "

Set-TextValue $ws.Range("D37") "485.86"
$ws.Range("E37").Value = "  -8.73%  "

$ws.Range("D38").Value = "3.132.73"
$ws.Range("E38").Value = "  +2.46%  "

$ws.Range("E39").Value = "  -4.61%  "

$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("E41").Value = "  -3.96%  "

Set-TextValue $ws.Range("D42") "8.08"
$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("E43").Value = "  -1.62%  "

Set-TextValue $ws.Range("D44") "0.252"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0533"
$ws.Range("E46").Value = "  +7.75%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D47") "121.22"
$ws.Range("E47").Value = "  +0.74%  "

Set-TextValue $ws.Range("D48") "2.01"
$ws.Range("E48").Value = "  -2.48%  "

Set-TextValue $ws.Range("D49") "24.36"
$ws.Range("E49").Value = "  +2.51%  "

$ws.Range("E50").Value = "  +2.03%  "

$ws.Range("E51").Value = "  -1.38%  "
